# Generate Report for Handoff
# The source file "f4932e3c-da9c-477e-9f19-8b9adeb55c80.md" has moved from
# "Handed back: in sync with en-US" to "Ready for handoff", with updated
# handoff timestamps, across the Overview sheet and each per-locale sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for f4932e3c-da9c-477e-9f19-8b9adeb55c80.md (row 3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-21 08:42:25"

# --- zh-cn sheet: row for f4932e3c-da9c-477e-9f19-8b9adeb55c80.md (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-21 08:42:22"

# --- de-de sheet: row for f4932e3c-da9c-477e-9f19-8b9adeb55c80.md (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-21 08:42:25"
